$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New match rows appended by the 03-01-2024 20:45 scraper run.
# Row 62 (Indice 61): Al Naser vs Khaitan
# Row 63 (Indice 62): Al Kuwait vs Al Shabab
# ---------------------------------------------------------------------

$rows = @(
    @{
        Row = 62
        Indice = 61
        Pais = "kuwait"
        Torneio = "premier-league"
        Temporada = "2023-2024"
        DataPartida = 45294.64930555555
        Home = "Al Naser"
        HomeGols = 1
        Away = "Khaitan"
        AwayGols = 1
        HomeOpenOdds = 1.35
        HomeOpenData = "03/01/2024 03:42"
        HomeCloseOdds = 1.37
        HomeCloseData = "03/01/2024 15:31"
        DrawOpenOdds = 4.9
        DrawOpenData = "03/01/2024 03:42"
        DrawCloseOdds = 5.28
        DrawCloseData = "03/01/2024 15:31"
        AwayOpenOdds = 6.72
        AwayOpenData = "03/01/2024 03:42"
        AwayCloseOdds = 6.02
        AwayCloseData = "03/01/2024 15:31"
        Url = "https://www.betexplorer.com/football/kuwait/premier-league/al-naser-khaitan/fea0Sbwh/"
    },
    @{
        Row = 63
        Indice = 62
        Pais = "kuwait"
        Torneio = "premier-league"
        Temporada = "2023-2024"
        DataPartida = 45294.75
        Home = "Al Kuwait"
        HomeGols = 6
        Away = "Al Shabab"
        AwayGols = 0
        HomeOpenOdds = 1.18
        HomeOpenData = "03/01/2024 06:12"
        HomeCloseOdds = 1.16
        HomeCloseData = "03/01/2024 17:52"
        DrawOpenOdds = 6.68
        DrawOpenData = "03/01/2024 06:12"
        DrawCloseOdds = 7.07
        DrawCloseData = "03/01/2024 17:52"
        AwayOpenOdds = 10.72
        AwayOpenData = "03/01/2024 06:12"
        AwayCloseOdds = 12
        AwayCloseData = "03/01/2024 17:52"
        Url = "https://www.betexplorer.com/football/kuwait/premier-league/al-kuwait-al-shabab/2mcF4eoo/"
    }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $srcRow = $rowNum - 1

    # Clone the formatting of the row directly above (same column styles:
    # bold/bordered/centered "Indice" cell in A, datetime-formatted cell in E)
    # so the new rows keep the same shared cellXfs instead of minting new ones.
    $ws.Range("A$srcRow").Copy() | Out-Null
    $ws.Range("A$rowNum").PasteSpecial(-4122) | Out-Null
    $ws.Range("E$srcRow").Copy() | Out-Null
    $ws.Range("E$rowNum").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($rowNum, 1).Value = $r.Indice
    $ws.Cells.Item($rowNum, 2).Value = $r.Pais
    $ws.Cells.Item($rowNum, 3).Value = $r.Torneio
    $ws.Cells.Item($rowNum, 4).Value = $r.Temporada
    $ws.Cells.Item($rowNum, 5).Value = $r.DataPartida
    $ws.Cells.Item($rowNum, 6).Value = $r.Home
    $ws.Cells.Item($rowNum, 7).Value = $r.HomeGols
    $ws.Cells.Item($rowNum, 8).Value = $r.Away
    $ws.Cells.Item($rowNum, 9).Value = $r.AwayGols
    $ws.Cells.Item($rowNum, 10).Value = $r.HomeOpenOdds
    $ws.Cells.Item($rowNum, 11).Value = $r.HomeOpenData
    $ws.Cells.Item($rowNum, 12).Value = $r.HomeCloseOdds
    $ws.Cells.Item($rowNum, 13).Value = $r.HomeCloseData
    $ws.Cells.Item($rowNum, 14).Value = $r.DrawOpenOdds
    $ws.Cells.Item($rowNum, 15).Value = $r.DrawOpenData
    $ws.Cells.Item($rowNum, 16).Value = $r.DrawCloseOdds
    $ws.Cells.Item($rowNum, 17).Value = $r.DrawCloseData
    $ws.Cells.Item($rowNum, 18).Value = $r.AwayOpenOdds
    $ws.Cells.Item($rowNum, 19).Value = $r.AwayOpenData
    $ws.Cells.Item($rowNum, 20).Value = $r.AwayCloseOdds
    $ws.Cells.Item($rowNum, 21).Value = $r.AwayCloseData
    $ws.Cells.Item($rowNum, 22).Value = $r.Url
}
